# purchase_list.xlsx update:
#  - Rewrite the "Ultrasonic Range Finder" line into two separate lines
#    (wide beam / narrow beam), re-price several items, and append new
#    rows for the Teensy board, an ODROID XU4, and a batch of prototyping
#    supplies (headers, breadboard, jumpers, wire kit).
#  - Move the Total row down to make room, widen column A, add a new
#    blank "Sheet2" tab (for the Ellison time card).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- row 5: Ultrasonic Range Finder (wide beam) ---------------------------
$ws.Range("A5").Value   = "Ultrasonic Range Finder (wide beam)"
$ws.Range("B5").Value   = "Obstacle avoidance prototyping"
$ws.Range("C5").Value   = 49.95
$ws.Range("D5").Value   = "10-20(?)"
$ws.Range("E5").Value   = 1
$ws.Range("F5").Formula = "=C5*E5"
$ws.Range("G5").Value   = "https://www.sparkfun.com/products/9491"

# --- row 6: Ultrasonic Range Finder (narrow beam) --------------------------
$ws.Range("A6").Value   = "Ultrasonic Range Finder (narrow beam)"
$ws.Range("B6").Value   = "Obstacle avoidance prototyping"
$ws.Range("C6").Value   = 49.95
$ws.Range("D6").Value   = "10-20(?)"
$ws.Range("E6").Value   = 1
$ws.Range("F6").Formula = "=C6*E6"
$ws.Range("G6").Value   = "https://www.sparkfun.com/products/9495"

# --- row 7: Teensy board ---------------------------------------------------
$ws.Range("A7").Value   = "Teensy board"
$ws.Range("B7").Value   = "Range finder computer interfacing"
$ws.Range("C7").Value   = 19.95
$ws.Range("D7").Value   = "10-30(?)"
$ws.Range("E7").Value   = 2
$ws.Range("F7").Formula = "=C7*E7"
$ws.Range("G7").Value   = "https://www.sparkfun.com/products/12646"

# --- row 8: ODROID XU4 ------------------------------------------------------
$ws.Range("A8").Value   = "ODROID XU4"
$ws.Range("B8").Value   = "Controls processing"
$ws.Range("C8").Value   = 75
$ws.Range("D8").Value   = 30
$ws.Range("E8").Value   = 1
$ws.Range("F8").Formula = "=C8*E8"
$ws.Range("G8").Value   = "http://www.hardkernel.com/main/products/prdt_info.php?g_code=G143452239825"

# --- rows 9-12: prototyping supplies (hyperlinked URLs, currency pricing) -
$ws.Range("A9").Value   = "Headers"
$ws.Range("B9").Value   = "Prototyping"
$ws.Range("C9").Value   = 1.5
$ws.Range("E9").Value   = 5
$ws.Range("F9").Formula = "=C9*E9"
$ws.Range("G9").Value   = "https://www.sparkfun.com/products/116"

$ws.Range("A10").Value   = "Breadboard"
$ws.Range("B10").Value   = "Prototyping"
$ws.Range("C10").Value   = 5.95
$ws.Range("E10").Value   = 1
$ws.Range("F10").Formula = "=C10*E10"
$ws.Range("G10").Value   = "https://www.sparkfun.com/products/12615"

$ws.Range("A11").Value   = "Jumpers"
$ws.Range("B11").Value   = "Prototyping"
$ws.Range("C11").Value   = 3.95
$ws.Range("E11").Value   = 1
$ws.Range("F11").Formula = "=C11*E11"
$ws.Range("G11").Value   = "https://www.sparkfun.com/products/9140"

$ws.Range("A12").Value   = "Wire Kit"
$ws.Range("B12").Value   = "Prototyping"
$ws.Range("C12").Value   = 6.95
$ws.Range("E12").Value   = 1
$ws.Range("F12").Formula = "=C12*E12"
$ws.Range("G12").Value   = "https://www.sparkfun.com/products/124"

# currency/accounting format + hyperlinks for the prototyping rows
$ws.Range("C9:C12,F9:F12").Style = "Currency"
$ws.Hyperlinks.Add($ws.Range("G9"),  "https://www.sparkfun.com/products/116")
$ws.Hyperlinks.Add($ws.Range("G10"), "https://www.sparkfun.com/products/12615")
$ws.Hyperlinks.Add($ws.Range("G11"), "https://www.sparkfun.com/products/9140")
$ws.Hyperlinks.Add($ws.Range("G12"), "https://www.sparkfun.com/products/124")

# --- blank spacer rows 13-17, Total row moves from 17 to 18 ---------------
$ws.Range("A13:G17").Value = ""
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""

$ws.Range("E18").Value   = "Total"
$ws.Range("F18").Formula = "=SUM(F2:F17)"

# --- column width / view tweaks -------------------------------------------
$ws.Columns("A").ColumnWidth = 33.1640625
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("F6").Select

# --- add the new (blank) Sheet2 for the Ellison time card -----------------
$wb.Worksheets.Add($null, $ws) | Out-Null
